# Bill of Materials - add "Recu" (received) tracking column (L) and a
# notes column (M); also correct the LED manufacturer/part/supplier-ref
# info for the green/red LEDs (rows 5 & 6) now that 3 ADC channels are
# chained via DMA and the LED refs needed fixing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RecuHeaderCell($addr) {
    $c = $ws.Range($addr)
    $c.Value = "Recu"
    $c.Font.Name = "Segoe UI"
    $c.Font.Size = 8
    $c.Font.Color = 0
    $c.Interior.ColorIndex = 15
    $c.HorizontalAlignment = -4108
    $c.Borders(7).LineStyle = 1
    $c.Borders(7).Weight = 2
    $c.Borders(10).LineStyle = 1
    $c.Borders(10).Weight = 2
}

function Set-RecuMarkCell($addr) {
    $c = $ws.Range($addr)
    $c.Value = "x"
    $c.Font.Name = "Segoe UI"
    $c.Font.Size = 8
    $c.Font.Color = 0
    $c.Borders(7).LineStyle = 1
    $c.Borders(7).Weight = 2
    $c.Borders(10).LineStyle = 1
    $c.Borders(10).Weight = 2
}

function Set-NoteCell($addr, $text) {
    $c = $ws.Range($addr)
    $c.Value = $text
    $c.Font.Name = "Segoe UI"
    $c.Font.Size = 8
    $c.Font.Color = 0
    $c.Borders(7).LineStyle = 1
    $c.Borders(7).Weight = 2
}

# --- New "Recu" column (L) ---------------------------------------------
Set-RecuHeaderCell "L1"
Set-RecuMarkCell "L2"
Set-RecuMarkCell "L4"
Set-RecuMarkCell "L7"
Set-RecuMarkCell "L8"
Set-RecuMarkCell "L9"
Set-RecuMarkCell "L10"
Set-RecuMarkCell "L11"
Set-RecuMarkCell "L12"
Set-RecuMarkCell "L13"
Set-RecuMarkCell "L15"

# --- Fix LED manufacturer / part number / supplier ref (3-channel ADC
#     chained DMA rework meant the LED refs needed a correction) --------
$ws.Range("G5").Value = "'KINGBRIGHT"
$ws.Range("I5").Value = "'KP-1608MGC"
$ws.Range("J5").Value = "'8529825"

$ws.Range("G6").Value = "'KINGBRIGHT"
$ws.Range("I6").Value = "'KP-1608SURCK"
$ws.Range("J6").Value = "'2290329"

# --- New notes column (M) ----------------------------------------------
Set-NoteCell "M3" "won't be ordered"
Set-NoteCell "M5" "ref issue!"
Set-NoteCell "M6" "ref issue!"

# --- Selection / active cell -------------------------------------------
$ws.Range("M4").Select()
